$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.119.06'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.285.23'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '535.02'
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").Value = '131.37'
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +3.68%  '
$ws.Range("D9").Value = '2.285.57'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").Value = '23.52'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = '2.690.23'
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").Value = '58.053.98'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '2.298.78'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '10.52'
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").Value = '313.03'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '6.45'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '63.11'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").Value = '7.96'
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("D28").Value = '1.27'
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("D29").Value = '170.69'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").Value = '1.07'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '288.13'
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '139.83'
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").Value = '0.0952'
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").Value = '18.09'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").Value = '10.94'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("E51").Value = '  +1.73%  '
